$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("종목명", "수요예측일", "희망공모가(원)", "확정공모가", "공모금액(백만)", "주간사"),
    @("케이웨더", "2024.01.12~01.18", "4,800~5,800", "-", 4800, "NH투자증권"),
    @("코셈", "2024.01.12~01.18", "12,000~14,000", "-", 7200, "키움증권"),
    @("이닉스(구.이닉스정호)", "2024.01.11~01.17", "9,200~11,000", "-", 27600, "삼성증권"),
    @("대신밸런스스팩17호", "2024.01.09~01.10", "2,000~2,000", "-", 11000, "대신증권"),
    @("에이치비인베스트먼트", "2024.01.08~01.12", "2,400~2,800", "-", 16000, "NH투자증권"),
    @("현대힘스", "2024.01.08~01.12", "5,000~6,300", "-", 43535, "미래에셋증권"),
    @("우진엔텍", "2024.01.08~01.12", "4,300~4,900", "-", 8858, "케이비증권"),
    @("포스뱅크", "2024.01.05~01.11", "13,000~15,000", "-", 19500, "하나증권"),
    @("스튜디오삼익", "2024.01.05~01.11", "14,500~16,500", "-", 12325, "DB금융투자"),
    @("하나스팩30호", "2023.12.08~12.11", "2,000~2,000", "2000", 14000, "하나증권"),
    @("디에스단석(구,단석산업)", "2023.12.05~12.11", "79,000~89,000", "100000", 96380, "KB증권,NH투자증권"),
    @("IBKS스팩23호", "2023.12.04~12.08", "2,000~2,000", "2000", 8000, "아이비케이투자증권"),
    @("블루엠텍", "2023.11.22~11.28", "15,000~19,000", "19000", 21000, "하나증권,키움증권"),
    @("LS머트리얼즈", "2023.11.22~11.28", "4,400~5,500", "6000", 64350, "키움증권,KB증권,이베스트투자증권,하이투자증권,NH투자증권"),
    @("삼성스팩9호", "2023.11.20~11.21", "2,000~2,000", "2000", 20000, "삼성증권"),
    @("교보스팩15호", "2023.11.20~11.21", "2,000~2,000", "2000", 7000, "교보증권"),
    @("케이엔에스", "2023.11.16~11.22", "19,000~22,000", "23000", 14250, "신영증권"),
    @("NH스팩30호", "2023.11.15~11.16", "2,000~2,000", "2000", 16000, "NH투자증권"),
    @("와이바이오로직스", "2023.11.10~11.16", "9,000~11,000", "9000", 13500, "유안타증권"),
    @("에이텀", "2023.11.09~11.15", "23,000~30,000", "18000", 14950, "하나증권")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $val = $data[$r][$c]
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        if ($val -is [string] -and $val -match "^[0-9]+$") {
            # Purely-numeric text (e.g. "2000") must stay a text cell,
            # not be auto-coerced into a number. Force text format for
            # the assignment, then strip the formatting it leaves behind
            # so the cell keeps its original (default) style.
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.ClearFormats()
        } else {
            $cell.Value = $val
        }
    }
}
